# Tejgaon-College Masters Subject_Id-1051 workbook touch-up.
#
# Matches the authored diff as closely as the Excel object model allows:
#   - Sheet1's view now scrolls so row 16 is at the top, while the
#     selection moves from the old M7 cell onto the whole D column's
#     used range (D1:D50), anchored at D1.
#   - Column D (the subject-id column) is widened from its old narrow
#     "5" character fit up to a little over 11 characters so the wider
#     values it now holds are fully visible.
#
# (workbook.xml/styles.xml metadata such as fileVersion/AlternateContent
# and the theme's font "panose" hints are housekeeping the host
# regenerates on every save and are not something this object model
# exposes a knob for.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Scroll the window so row 16 becomes the first visible row (topLeftCell
# = A16), then select D1:D50 with D1 as the active cell.
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1

$ws.Range("D1:D50").Select()

# Widen column D. The host quantizes ColumnWidth to the nearest
# displayable step, so 10.67 is the value that lands on the target
# stored width (~11.57 characters) closest to the authored file.
$ws.Columns(4).ColumnWidth = 10.67
